$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the two empty helper sheets (Planilha4, Planilha5)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Planilha4").Delete()
$wb.Worksheets.Item("Planilha5").Delete()

# ---------------------------------------------------------------------------
# 2. Planilha3: insert a new header row above the existing grid and fill in
#    the new "maze" markers / values
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Planilha3")

# Push everything down one row (keeps styles of rows 5..17 intact, they just
# become rows 6..18)
$ws.Rows("5:5").Insert()

# Row 5 had no style of its own after the insert - borrow the look of row 18
# (same column layout: E/F=26, G..K=27, L/M=28)
$ws.Range("E18:M18").Copy()
$ws.Range("E5:M5").PasteSpecial(-4122)

# New header values for row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 3

# Updated grid values / new "#" wall markers (rows 7-17 are the shifted
# original grid, rows were previously 6-16)
$ws.Range("K7").Value = "#"
$ws.Range("L7").Value = "#"
$ws.Range("M7").Value = "#"
$ws.Range("N7").Value = "#"

$ws.Range("B8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("J8").Value = "#"
$ws.Range("N8").Value = "#"

$ws.Range("B9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("I9").Value = "#"
$ws.Range("N9").Value = "#"

$ws.Range("B10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("H10").Value = "#"
$ws.Range("N10").Value = "#"

$ws.Range("B11").Value = 3
$ws.Range("G11").Value = "#"
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = 3
$ws.Range("N11").Value = "#"

$ws.Range("B12").Value = 3
$ws.Range("F12").Value = "#"
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = 3
$ws.Range("N12").Value = "#"

$ws.Range("B13").Value = 2
$ws.Range("F13").Value = "#"
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 2
$ws.Range("K13").Value = 2
$ws.Range("N13").Value = "#"

$ws.Range("B14").Value = 2
$ws.Range("F14").Value = "#"
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 2
$ws.Range("N14").Value = "#"

$ws.Range("B15").Value = 2
$ws.Range("F15").Value = "#"
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 2
$ws.Range("M15").Value = "#"

$ws.Range("B16").Value = 2
$ws.Range("E16").Value = "#"
$ws.Range("F16").Value = "#"
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = "#"

$ws.Range("D17").Value = "#"
$ws.Range("F17").Value = "#"
$ws.Range("G17").Value = "#"
$ws.Range("H17").Value = "#"
$ws.Range("I17").Value = "#"
$ws.Range("J17").Value = "#"
$ws.Range("K17").Value = "#"

# New small legend table below the grid
$ws.Range("E21").Value = "diag"
$ws.Range("F21").Value = "up"
$ws.Range("G21").Value = "down"
$ws.Range("H21").Value = "right"

$ws.Range("E22").Value = 4.5
$ws.Range("F22").Value = 5.5
$ws.Range("G22").Value = 5.5

# Restore the view: scroll/selection on this sheet
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F12").Select()

# ---------------------------------------------------------------------------
# 3. Planilha2: a handful of "#" style xfs were merged/cleaned up - column AC
#    (the mirrored "V3" column) loses its dedicated numeric-format style and
#    reverts to the default style, while AC28/AC54 (blank marker cells) move
#    onto the freed-up style slot.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Planilha2")
for ($r = 3; $r -le 78; $r++) {
    if ($r -eq 27 -or $r -eq 53) { continue }
    $cell = $ws2.Cells.Item($r, 29)
    if ($cell.Value2 -ne $null) {
        $cell.Style = "Normal"
    }
}

# Make the active sheet/tab match the saved workbook state (Planilha3 is now
# the last, 4th, sheet and should be the one shown on open)
$ws.Activate()
